$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report contains several duplicated stock-report line items (same
# product code family) whose rows got their data fields crossed when the
# report was regenerated. For each group of rows below, the B (item code),
# C (description), D (rate), E (rate2), F (qty) and G (value) fields need
# to be cyclically rotated one position among the rows in the group
# (row[i] <- row[i-1], wrapping around) to restore the correct pairing.
$cols = "B", "C", "D", "E", "F", "G"

$g1 = @(149,150)
$g2 = @(313,314)
$g3 = @(316,317)
$g4 = @(350,351,352)
$g5 = @(355,356)
$g6 = @(372,373)
$g7 = @(389,390)
$g8 = @(419,420)
$g9 = @(421,422)
$g10 = @(431,432)
$g11 = @(457,458)
$g12 = @(579,580)
$g13 = @(583,584)
$g14 = @(586,587)
$g15 = @(590,591)
$g16 = @(593,594)
$g17 = @(601,602)
$g18 = @(687,688)
$g19 = @(709,710)
$g20 = @(715,716)
$g21 = @(720,721)
$g22 = @(859,860)

$rowGroups = @($g1, $g2, $g3, $g4, $g5, $g6, $g7, $g8, $g9, $g10, $g11, $g12, $g13, $g14, $g15, $g16, $g17, $g18, $g19, $g20, $g21, $g22)

foreach ($rows in $rowGroups) {
    $n = $rows.Count
    foreach ($col in $cols) {
        $vals = @()
        foreach ($r in $rows) { $vals += ,($ws.Range("$col$r").Value2) }
        for ($i = 0; $i -lt $n; $i++) {
            $srcIdx = ($i - 1 + $n) % $n
            $ws.Range("$col$($rows[$i])").Value2 = $vals[$srcIdx]
        }
    }
}
